# "#5: cash & deposit done"
#
# Sheet4 (現金 / cash) and Sheet5 (存款 / deposit) previously only stored a
# few columns and their header row (row 1) simply echoed the first data
# row. This change brings both sheets in line with the other property
# sheets (land/building/car): the header row now holds the real field
# names, and every data row gets the full set of trailing metadata
# columns (property_category/category/date/legislator_name/
# legislator_id/source_file/index). Sheet5's F4 value also switches from
# being stored as text to a genuine number.

$wb = $excel.ActiveWorkbook

# xlPasteFormats = -4122 ; used to copy a cell's number format/style onto
# a newly written cell without touching its value.
$xlPasteFormats = -4122

# -----------------------------------------------------------------
# Sheet4 : 現金 (cash)
# -----------------------------------------------------------------
$ws4 = $wb.Worksheets.Item(4)

# --- header row (row 1) ---
$ws4.Cells.Item(1, 2).Value = "currency"
$ws4.Cells.Item(1, 3).Value = "owner"
$ws4.Cells.Item(1, 4).Value = "total"

$ws4.Cells.Item(1, 5).Value = "property_category"
$ws4.Cells.Item(1, 6).Value = "category"
$ws4.Cells.Item(1, 7).Value = "date"
$ws4.Cells.Item(1, 8).Value = "legislator_name"
$ws4.Cells.Item(1, 9).Value = "legislator_id"
$ws4.Cells.Item(1, 10).Value = "source_file"
$ws4.Cells.Item(1, 11).Value = "index"

# give the newly-created header cells (E1:K1) the same style as the
# existing header cells (bold, bordered -> style index 1)
$ws4.Cells.Item(1, 2).Copy()
$ws4.Range($ws4.Cells.Item(1, 5), $ws4.Cells.Item(1, 11)).PasteSpecial($xlPasteFormats)

# --- data row (row 2) ---
# B2/C2/D2 already hold the correct values (新臺幣 / 林淑芬 / 2320000)
$ws4.Cells.Item(2, 5).Value = "cash"
$ws4.Cells.Item(2, 6).Value = "normal"

# "2013-12-30" must stay a literal text value (like the other sheets),
# not get auto-converted into an Excel date serial number.
$ws4.Cells.Item(2, 7).NumberFormat = "@"
$ws4.Cells.Item(2, 7).Value = "2013-12-30"

$ws4.Cells.Item(2, 8).Value = "林淑芬"
$ws4.Cells.Item(2, 9).Value = 1337
$ws4.Cells.Item(2, 10).Value = "tmp63cf1"
$ws4.Cells.Item(2, 11).Value = 41

# give the newly-created data cells (E2:K2) the same style as the
# existing data cells (style index 2) - this also clears the temporary
# "@" text format applied to G2 above.
$ws4.Cells.Item(2, 2).Copy()
$ws4.Range($ws4.Cells.Item(2, 5), $ws4.Cells.Item(2, 11)).PasteSpecial($xlPasteFormats)

# -----------------------------------------------------------------
# Sheet5 : 存款 (deposit)
# -----------------------------------------------------------------
$ws5 = $wb.Worksheets.Item(5)

# --- header row (row 1) : switch from echoed data to real field names ---
$ws5.Cells.Item(1, 2).Value = "bank"
$ws5.Cells.Item(1, 3).Value = "deposit_type"
$ws5.Cells.Item(1, 4).Value = "currency"
$ws5.Cells.Item(1, 5).Value = "owner"
$ws5.Cells.Item(1, 6).Value = "total"

$ws5.Cells.Item(1, 7).Value = "property_category"
$ws5.Cells.Item(1, 8).Value = "category"
$ws5.Cells.Item(1, 9).Value = "date"
$ws5.Cells.Item(1, 10).Value = "legislator_name"
$ws5.Cells.Item(1, 11).Value = "legislator_id"
$ws5.Cells.Item(1, 12).Value = "source_file"
$ws5.Cells.Item(1, 13).Value = "index"

$ws5.Cells.Item(1, 2).Copy()
$ws5.Range($ws5.Cells.Item(1, 7), $ws5.Cells.Item(1, 13)).PasteSpecial($xlPasteFormats)

# --- row 2 (index 46) : B2:F2 keep their current values ---
$ws5.Cells.Item(2, 7).Value = "deposit"
$ws5.Cells.Item(2, 8).Value = "normal"

$ws5.Cells.Item(2, 9).NumberFormat = "@"
$ws5.Cells.Item(2, 9).Value = "2013-12-30"

$ws5.Cells.Item(2, 10).Value = "林淑芬"
$ws5.Cells.Item(2, 11).Value = 1337
$ws5.Cells.Item(2, 12).Value = "tmp63cf1"
$ws5.Cells.Item(2, 13).Value = 46

$ws5.Cells.Item(2, 2).Copy()
$ws5.Range($ws5.Cells.Item(2, 7), $ws5.Cells.Item(2, 13)).PasteSpecial($xlPasteFormats)

# --- row 3 (index 47) : B3:F3 keep their current values ---
$ws5.Cells.Item(3, 7).Value = "deposit"
$ws5.Cells.Item(3, 8).Value = "normal"

$ws5.Cells.Item(3, 9).NumberFormat = "@"
$ws5.Cells.Item(3, 9).Value = "2013-12-30"

$ws5.Cells.Item(3, 10).Value = "林淑芬"
$ws5.Cells.Item(3, 11).Value = 1337
$ws5.Cells.Item(3, 12).Value = "tmp63cf1"
$ws5.Cells.Item(3, 13).Value = 47

$ws5.Cells.Item(3, 2).Copy()
$ws5.Range($ws5.Cells.Item(3, 7), $ws5.Cells.Item(3, 13)).PasteSpecial($xlPasteFormats)

# --- row 4 (index 48) : F4 used to be stored as text "1420000";
#     it becomes a genuine number, B4:E4 stay as-is ---
$ws5.Cells.Item(4, 6).Value = 1420000

$ws5.Cells.Item(4, 7).Value = "deposit"
$ws5.Cells.Item(4, 8).Value = "normal"

$ws5.Cells.Item(4, 9).NumberFormat = "@"
$ws5.Cells.Item(4, 9).Value = "2013-12-30"

$ws5.Cells.Item(4, 10).Value = "林淑芬"
$ws5.Cells.Item(4, 11).Value = 1337
$ws5.Cells.Item(4, 12).Value = "tmp63cf1"
$ws5.Cells.Item(4, 13).Value = 48

$ws5.Cells.Item(4, 2).Copy()
$ws5.Range($ws5.Cells.Item(4, 7), $ws5.Cells.Item(4, 13)).PasteSpecial($xlPasteFormats)
